$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet 1: LP1912
# -----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Cells.Item(2,1).Value = "Última actualización: 03:29:38"
$ws1.Cells.Item(3,1).Value = "Total filas: 7"

# Row 6
$ws1.Cells.Item(6,1).Value = "03:29:38"
$ws1.Cells.Item(6,2).Value = "03:31"
$ws1.Cells.Item(6,3).Value = "14_ABASTO"
$ws1.Cells.Item(6,4).Value = 2
$ws1.Cells.Item(6,5).Value = "LP1912"

# Row 7
$ws1.Cells.Item(7,1).Value = "03:29:38"
$ws1.Cells.Item(7,2).Value = "04:01"
$ws1.Cells.Item(7,3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(7,4).Value = 32
$ws1.Cells.Item(7,5).Value = "LP1912"

# Row 8
$ws1.Cells.Item(8,1).Value = "03:29:38"
$ws1.Cells.Item(8,2).Value = "04:45"
$ws1.Cells.Item(8,3).Value = "215A_EL PATO"
$ws1.Cells.Item(8,4).Value = 76
$ws1.Cells.Item(8,5).Value = "LP1912"

# Row 9
$ws1.Cells.Item(9,1).Value = "03:29:38"
$ws1.Cells.Item(9,2).Value = "04:53"
$ws1.Cells.Item(9,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(9,4).Value = 84
$ws1.Cells.Item(9,5).Value = "LP1912"

# Row 10
$ws1.Cells.Item(10,1).Value = "03:29:38"
$ws1.Cells.Item(10,2).Value = "05:16"
$ws1.Cells.Item(10,3).Value = "17_ROMERO"
$ws1.Cells.Item(10,4).Value = 107
$ws1.Cells.Item(10,5).Value = "LP1912"

# Row 11 (new)
$ws1.Cells.Item(11,1).Value = "03:29:38"
$ws1.Cells.Item(11,2).Value = "05:21"
$ws1.Cells.Item(11,3).Value = "14_ABASTO"
$ws1.Cells.Item(11,4).Value = 112
$ws1.Cells.Item(11,5).Value = "LP1912"

# Row 12 (new)
$ws1.Cells.Item(12,1).Value = "03:29:38"
$ws1.Cells.Item(12,2).Value = "05:22"
$ws1.Cells.Item(12,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(12,4).Value = 113
$ws1.Cells.Item(12,5).Value = "LP1912"

# -----------------------------------------------------------------
# Sheet 2: LP1912-215
# -----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Cells.Item(2,1).Value = "Última actualización: 03:29:38"
$ws2.Cells.Item(3,1).Value = "Total filas: 1"

# Row 6 becomes what used to be row 7 (215A_EL PATO / 04:45 / 76)
$ws2.Cells.Item(6,1).Value = "03:29:38"
$ws2.Cells.Item(6,2).Value = "04:45"
$ws2.Cells.Item(6,3).Value = "215A_EL PATO"
$ws2.Cells.Item(6,4).Value = 76
$ws2.Cells.Item(6,5).Value = "LP1912"

# Remove the old row 7 entirely
$ws2.Rows(7).Delete()

# -----------------------------------------------------------------
# Sheet 3: 6203-6173
# -----------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Cells.Item(2,1).Value = "Última actualización: 03:29:38"
